$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "wide" layout entirely before writing the new "long" layout.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Years"
$ws.Range("B1").Value = "Non_Hispanic "
$ws.Range("C1").Value = "Hispanic "

# Data rows: Years, Non_Hispanic, Hispanic
$ws.Range("A2").Value = 2015
$ws.Range("B2").Value = 452140
$ws.Range("C2").Value = 112568

$ws.Range("A3").Value = 2016
$ws.Range("B3").Value = 428629
$ws.Range("C3").Value = 121299

$ws.Range("A4").Value = 2017
$ws.Range("B4").Value = 432634
$ws.Range("C4").Value = 118362

$ws.Range("A5").Value = 2018
$ws.Range("B5").Value = 430354
$ws.Range("C5").Value = 122476

$ws.Range("A6").Value = 2019
$ws.Range("B6").Value = 443100
$ws.Range("C6").Value = 124615

$ws.Range("A7").Value = 2020
$ws.Range("B7").Value = 450107
$ws.Range("C7").Value = 130348

# Apply left alignment to the whole used range (new style xf).
$ws.Range("A1:C7").HorizontalAlignment = -4131

# Column widths (target stored widths are 18.5703125 / 17.140625; the
# ColumnWidth setter here snaps to 1/6-character increments, so feed the
# input that lands on the closest achievable snapped value).
$ws.Columns.Item(2).ColumnWidth = 17.59
$ws.Columns.Item(3).ColumnWidth = 16.25

# Selection
$ws.Range("B1").Select()
